$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.6579528450965881
$ws.Range("B1").Value = 1.493727684020996
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.10013484954834
$ws.Range("E1").Value = 1.294928789138794
